$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures.
# D-column price values are forced as text via a leading apostrophe
# (matches the workbook convention of storing these as inline strings,
# not numbers - some values like "26.901.14" cannot be numbers anyway,
# and this keeps all D-column entries consistently typed as text).
$ws.Range("D2").Value = "'26.901.14"
$ws.Range("D3").Value = "'1.549.37"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D5").Value = "'206.51"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'0.486"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'22.15"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "'1.770.84"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'1.549.64"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'26.906.53"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'217.33"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'0.0₃0697"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "'154.24"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'14.91"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "'1.415.23"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").Value = "'0.965"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").Value = "'0.525"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'5.70"
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").Value = "'64.45"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "'1.684.17"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'87.53"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  +0.09%  "
